$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17, shifting existing rows 17-61 down to 18-62.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with data.
$ws.Cells.Item(17, 1).Value = 4
$ws.Cells.Item(17, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(17, 3).Value = "Los Lagos"
$ws.Cells.Item(17, 4).Value = 44624
$ws.Cells.Item(17, 5).Value = 10
$ws.Cells.Item(17, 6).Value = 100112031
$ws.Cells.Item(17, 7).Value = "Poroto verde"
$ws.Cells.Item(17, 8).Value = "Magnum"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 40
$ws.Cells.Item(17, 11).Value = 30000
$ws.Cells.Item(17, 12).Value = 30000
$ws.Cells.Item(17, 13).Value = 30000
$ws.Cells.Item(17, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(17, 15).Value = "Región Metropolitana"
$ws.Cells.Item(17, 16).Value = 1200
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Apply the same style as other date cells (D column) to the new D17 cell.
$ws.Cells.Item(17, 4).NumberFormat = $ws.Cells.Item(18, 4).NumberFormat
